$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.220.72"
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("D3").Value = "2.165.02"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "236.04"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("E6").Value = "  -3.12%  "
$ws.Range("D7").Value = "69.36"
$ws.Range("E7").Value = "  -5.40%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -6.68%  "
$ws.Range("D10").Value = "38.75"
$ws.Range("E10").Value = "  -9.25%  "
$ws.Range("E11").Value = "  -4.84%  "
$ws.Range("D12").Value = "54.28"
$ws.Range("E12").Value = "  -5.53%  "
$ws.Range("D13").Value = "0.100"
$ws.Range("E13").Value = "  -2.99%  "
$ws.Range("D14").Value = "6.67"
$ws.Range("E14").Value = "  -5.80%  "
$ws.Range("D15").Value = "2.486.74"
$ws.Range("E15").Value = "  -2.40%  "
$ws.Range("D16").Value = "14.31"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "2.131.89"
$ws.Range("E17").Value = "  -3.76%  "
$ws.Range("D18").Value = "0.783"
$ws.Range("E18").Value = "  -6.44%  "
$ws.Range("D19").Value = "41.053.66"
$ws.Range("E19").Value = "  -1.95%  "
$ws.Range("D20").Value = "0.0₃0994"
$ws.Range("E20").Value = "  -7.49%  "
$ws.Range("D21").Value = "69.77"
$ws.Range("E21").Value = "  -4.28%  "
$ws.Range("D22").Value = "5.78"
$ws.Range("E22").Value = "  -6.17%  "
$ws.Range("D23").Value = "224.59"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").Value = "9.31"
$ws.Range("E24").Value = "  -14.85%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("E26").Value = "  -10.68%  "
$ws.Range("D27").Value = "10.63"
$ws.Range("E27").Value = "  -9.83%  "
$ws.Range("D28").Value = "3.30"
$ws.Range("E28").Value = "  -8.47%  "
$ws.Range("E29").Value = "  -3.57%  "
$ws.Range("D30").Value = "2.16"
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("D31").Value = "168.42"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("E32").Value = "  -3.52%  "
$ws.Range("D33").Value = "30.09"
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("E34").Value = "  -5.34%  "
$ws.Range("D35").Value = "5.08"
$ws.Range("E35").Value = "  -9.01%  "
$ws.Range("D36").Value = "0.120"
$ws.Range("E36").Value = "  -3.88%  "
$ws.Range("E37").Value = "  -6.70%  "
$ws.Range("D38").Value = "4.08"
$ws.Range("E38").Value = "  -4.28%  "
$ws.Range("D39").Value = "0.0279"
$ws.Range("E39").Value = "  -6.85%  "
$ws.Range("E40").Value = "  -3.28%  "
$ws.Range("D41").Value = "11.45"
$ws.Range("E41").Value = "  -16.25%  "
$ws.Range("D42").Value = "5.30"
$ws.Range("E42").Value = "  -5.68%  "
$ws.Range("D43").Value = "57.95"
$ws.Range("E43").Value = "  -11.92%  "
$ws.Range("E44").Value = "  -6.04%  "
$ws.Range("D45").Value = "8.22"
$ws.Range("E45").Value = "  -5.66%  "
$ws.Range("E46").Value = "  -4.52%  "
$ws.Range("D47").Value = "96.59"
$ws.Range("E47").Value = "  -7.59%  "
$ws.Range("E48").Value = "  -4.13%  "
$ws.Range("E49").Value = "  -4.55%  "
$ws.Range("D50").Value = "2.16"
$ws.Range("E50").Value = "  -9.65%  "
$ws.Range("D51").Value = "2.60"
$ws.Range("E51").Value = "  -3.36%  "
